$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Khoni")

# Add a new "2023" column (K) to the table, reusing the formatting of the
# preceding "2022" column (J).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 589.29999999999995
$ws.Range("K5").Value = 500.8
$ws.Range("K6").Value = 719.2
